# Update J1 header text from "Studentoffer" to "Specialoffer"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Specialoffer"

# Update the view state: scroll so column B is the leftmost visible column,
# and select I8 as the active cell.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("I8").Select()
